# Simplify user testing template: clean sheet names and titles
#  - Sheet tabs renamed from "Test N" to just "N"
#  - A1 title on each sheet changed from "AEGIS Performance Testing - Test #N"
#    to plain "AEGIS Performance Testing"

$wb = $excel.ActiveWorkbook

$count = $wb.Worksheets.Count
for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = "$i"
    $ws.Range("A1").Value = "AEGIS Performance Testing"
}
